$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet at the very front of the
#    workbook (before "ODI Batting" / "ODI Bowling").
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Style the header row to match the other header rows in the
# workbook (bold, thin box border, centered/top aligned).
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$playerInfo.Range("A2").Value = "'4511"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Scott M Boland"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and
#    replace the full scorecard URL with just the numeric match code.
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "3874"
    3 = "3875"
    4 = "3876"
    5 = "3878"
    6 = "3886"
    7 = "3888"
    8 = "3903"
    9 = "3905"
    10 = "3929"
    11 = "3931"
    12 = "3937"
    13 = "3940"
    14 = "3947"
    15 = "3950"
}

foreach ($row in $battingCodes.Keys) {
    $cell = $batting.Range("D$row")
    $cell.Value = "'" + $battingCodes[$row]
    $cell.Style = "Normal"
}

# ------------------------------------------------------------------
# 3. "ODI Bowling" sheet: same MATCH_CARD_LINK -> MATCH_CODE rename
#    and URL -> numeric code replacement (column B here).
# ------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "3874"
    3 = "3875"
    4 = "3876"
    5 = "3878"
    6 = "3886"
    7 = "3888"
    8 = "3905"
    9 = "3929"
    10 = "3931"
    11 = "3937"
    12 = "3940"
    13 = "3947"
    14 = "3950"
}

foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowling.Range("B$row")
    $cell.Value = "'" + $bowlingCodes[$row]
    $cell.Style = "Normal"
}
